# Applies the weekly Fruta/Hortaliza update to the Chirimoya sheet:
# existing rows 2-14 are refreshed with the latest field values and two
# new price records are appended as rows 15-16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target values for rows 2-16, columns A (1) .. T (20)
$data = @(
    @(5, 'Macroferia Regional de Talca', 'Maule', 44432, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 70, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44454, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 320, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44454, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 300, 28000, 28000, 28000, '$/bandeja 10 kilos', 'Provincia de Limarí', 2800, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44421, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 30, 35000, 35000, 35000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3500, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44431, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 30, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44435, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 160, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44445, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 250, 28000, 30000, 29200, '$/bandeja 10 kilos', 'Provincia de Limarí', 2920, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44446, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 200, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44434, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 60, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44441, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 150, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44453, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 135, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44448, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 100, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44448, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 80, 28000, 28000, 28000, '$/bandeja 10 kilos', 'Provincia de Limarí', 2800, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44447, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Especial', 50, 32000, 32000, 32000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3200, 10),
    @(5, 'Macroferia Regional de Talca', 'Maule', 44438, 7, 'Fruta', 100107, 'Otros', 100107002, 'Chirimoya', 'Cultivar IV Región', 'Primera', 100, 30000, 30000, 30000, '$/bandeja 10 kilos', 'Provincia de Limarí', 3000, 10),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
    # Column D (date) keeps the existing custom date/time number format
    $ws.Cells.Item($rowIndex, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
}

Write-Host 'Chirimoya sheet updated: rows 2-14 refreshed, rows 15-16 added.'
